# Apply the commit "Act greficos y tablas web pob":
#  - Rename sheets: "Datos" -> "Data", "Ficha técnica" -> "Metadata"
#  - Reverse (mirror) the year ordering + associated data on the Data sheet
#  - Rework the Metadata sheet field names to lowercase keys, split
#    "CITA" out into its own "observaciones"/"cita" rows, and append a
#    final attribution row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Datos" -> "Data"
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsData.Name = "Data"

# The year labels (column A, rows 2-14) get reversed in order, and the
# associated B/C values follow the same row they were attached to (i.e.
# the whole table is mirrored top-to-bottom around row 8).
# Years must remain text (shared-string) cells, like the originals, not
# numbers - so write them through a helper cell using a text formula and
# paste-as-values, which avoids Excel's "looks like a number" auto-typing
# and keeps the destination cell's formatting untouched.
$years = @("2019","2018","2017","2016","2015","2014","2013","2012","2011","2010","2009","2008","2007")
$helper = $wsData.Cells.Item(50, 1)
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 2 + $i
    $helper.Formula = '=TEXT(' + $years[$i] + ',"0")'
    $helper.Copy()
    $wsData.Cells.Item($row, 1).PasteSpecial(-4163)
}
$helper.Clear()

$bValues = @(10.8, 10.2, 11.3, 9.9, 12.8, 11.6, 14.3, 13.9, 12.7, 13.9, 14.1, 15.3, 13.7)
$cValues = @(3.7, 4, 3.9, 4.5, 4.2, 4.1, 3.7, 4.5, 4.2, 5.1, 4.5, 4.4, 4.5)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = 2 + $i
    $wsData.Cells.Item($row, 2).Value = $bValues[$i]
    $wsData.Cells.Item($row, 3).Value = $cValues[$i]
}

# ---------------------------------------------------------------
# Sheet 2: "Ficha técnica" -> "Metadata"
# ---------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item(2)
$wsMeta.Name = "Metadata"

$wsMeta.Cells.Item(2, 1).Value = "nomindicador"
$wsMeta.Cells.Item(2, 2).Value = "Porcentaje de personas que viven en asentamientos"

$wsMeta.Cells.Item(3, 1).Value = "derecho"
$wsMeta.Cells.Item(3, 2).Value = "Vivienda"

$wsMeta.Cells.Item(4, 1).Value = "conindicador"
$wsMeta.Cells.Item(4, 2).Value = "Asentamientos"

$wsMeta.Cells.Item(5, 1).Value = "tipoind"
$wsMeta.Cells.Item(5, 2).Value = "Resultados"

$wsMeta.Cells.Item(6, 1).Value = "definicion"
$wsMeta.Cells.Item(6, 2).Value = "El indicador mide el porcentaje de personas en viviendas ubicadas en asentamiento irregular."

$wsMeta.Cells.Item(7, 1).Value = "calculo"
$wsMeta.Cells.Item(7, 2).Value = "Para cada año calcular: (Cantidad de personas que residen en viviendas ubicadas en asentamiento irregular / Cantidad total de personas en viviendas particulares)*100"

$wsMeta.Cells.Item(8, 1).Value = "observaciones"
$wsMeta.Cells.Item(8, 2).Value = "Sin observaciones"

$wsMeta.Cells.Item(9, 1).Value = "cita"
$wsMeta.Cells.Item(9, 2).Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"

$wsMeta.Cells.Item(10, 1).Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$wsMeta.Cells.Item(10, 2).Value = " "
